# Update the label in A1 (merged A1:A2) from "Spatial smoother terms" to "Model components"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Model components"

# Move the sheet's active selection to the merged header cell A1:A2
$ws.Range("A1:A2").Select()
